# Auto-generated Excel COM-interop edit script
# Applies updated market/profit figures to the Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 402.8
$ws.Range("I11").Value = 402.8
$ws.Range("K11").Value = 402.8
$ws.Range("M11").Value = -262.8
$ws.Range("H33").Value = 3979.6667
$ws.Range("I33").Value = 3979.6667
$ws.Range("K33").Value = 3979.6667
$ws.Range("M33").Value = -3750.6667
$ws.Range("H70").Value = 528.2273
$ws.Range("I70").Value = 419.47058
$ws.Range("J70").Value = 898
$ws.Range("K70").Value = 1258.41174
$ws.Range("L70").Value = 2694
$ws.Range("M70").Value = -988.41174
$ws.Range("N70").Value = -3234
$ws.Range("H73").Value = 528.2273
$ws.Range("I73").Value = 419.47058
$ws.Range("J73").Value = 898
$ws.Range("K73").Value = 1258.41174
$ws.Range("L73").Value = 2694
$ws.Range("M73").Value = -322.41174
$ws.Range("N73").Value = -4566
$ws.Range("H80").Value = 390.75
$ws.Range("I80").Value = 221.16667
$ws.Range("J80").Value = 560.3333
$ws.Range("K80").Value = 663.50001
$ws.Range("L80").Value = 1680.9999
$ws.Range("M80").Value = 334.49999
$ws.Range("N80").Value = -3676.9999
$ws.Range("H83").Value = 390.75
$ws.Range("I83").Value = 221.16667
$ws.Range("J83").Value = 560.3333
$ws.Range("K83").Value = 1990.50003
$ws.Range("L83").Value = 5042.9997
$ws.Range("M83").Value = 3001.49997
$ws.Range("N83").Value = -15026.9997
$ws.Range("H112").Value = 4679.1875
$ws.Range("I112").Value = 4211
$ws.Range("J112").Value = 4746.0713
$ws.Range("K112").Value = 12633
$ws.Range("L112").Value = 14238.2139
$ws.Range("M112").Value = -11525
$ws.Range("N112").Value = -16454.2139
$ws.Range("H132").Value = 2708.7856
$ws.Range("I132").Value = 1285.4546
$ws.Range("J132").Value = 7927.6665
$ws.Range("K132").Value = 3856.3638
$ws.Range("L132").Value = 23782.9995
$ws.Range("M132").Value = -1326.3638
$ws.Range("N132").Value = -28842.9995
$ws.Range("H138").Value = 2144.1892
$ws.Range("I138").Value = 1447.2273
$ws.Range("K138").Value = 4341.6819
$ws.Range("M138").Value = 798.3181000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2287.5173
$ws.Range("I2").Value = 2077.2
$ws.Range("K2").Value = 2077.2
$ws.Range("M2").Value = -1964.2
$ws.Range("H32").Value = 28936.672
$ws.Range("I32").Value = 16648.791
$ws.Range("K32").Value = 16648.791
$ws.Range("M32").Value = -16361.791
$ws.Range("H102").Value = 2441.5293
$ws.Range("I102").Value = 2441.5293
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2441.5293
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = -819.5293000000001
$ws.Range("M102").ClearContents()
$ws.Range("H110").Value = 2359.225
$ws.Range("I110").Value = 2422.8057
$ws.Range("K110").Value = 2422.8057
$ws.Range("M110").Value = -377.8056999999999
$ws.Range("H116").Value = 2287.5173
$ws.Range("I116").Value = 2077.2
$ws.Range("K116").Value = 2077.2
$ws.Range("M116").Value = 216.8000000000002
$ws.Range("H139").Value = 109331.336
$ws.Range("J139").Value = 69999.5
$ws.Range("L139").Value = 69999.5
$ws.Range("N139").Value = -80279.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2287.5173
$ws.Range("I3").Value = 2077.2
$ws.Range("K3").Value = 2077.2
$ws.Range("M3").Value = -1963.2
$ws.Range("H94").Value = 1254.3636
$ws.Range("I94").Value = 1254.3636
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1254.3636
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = -803.3635999999999
$ws.Range("M94").ClearContents()
$ws.Range("H134").Value = 1360.625
$ws.Range("I134").Value = 1232.9
$ws.Range("K134").Value = 3698.7
$ws.Range("M134").Value = -1163.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1511.4783
$ws.Range("I16").Value = 1285.4375
$ws.Range("K16").Value = 1285.4375
$ws.Range("M16").Value = -998.4375
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = 150
$ws.Range("M22").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("N48").Value = 0
$ws.Range("L48").ClearContents()
$ws.Range("H58").Value = 1074.25
$ws.Range("I58").Value = 992.0909
$ws.Range("K58").Value = 992.0909
$ws.Range("M58").Value = -789.0909
$ws.Range("H107").Value = 1102.0938
$ws.Range("I107").Value = 1215.875
$ws.Range("K107").Value = 1215.875
$ws.Range("M107").Value = 704.125
$ws.Range("H113").Value = 1511.4783
$ws.Range("I113").Value = 1285.4375
$ws.Range("K113").Value = 1285.4375
$ws.Range("M113").Value = 884.5625
$ws.Range("H134").Value = 4559.846
$ws.Range("I134").Value = 4752.5454
$ws.Range("K134").Value = 14257.6362
$ws.Range("M134").Value = -11722.6362
$ws.Range("H136").Value = 1074.25
$ws.Range("I136").Value = 992.0909
$ws.Range("K136").Value = 2976.2727
$ws.Range("M136").Value = -426.2727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 285
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 285
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = 855
$ws.Range("N23").Value = -1325
$ws.Range("L23").ClearContents()
$ws.Range("H38").Value = 310
$ws.Range("I38").Value = 250
$ws.Range("J38").Value = 322
$ws.Range("K38").Value = 750
$ws.Range("L38").Value = 966
$ws.Range("M38").Value = -403
$ws.Range("N38").Value = -1660
$ws.Range("H98").Value = 1387.875
$ws.Range("I98").Value = 788.6
$ws.Range("K98").Value = 2365.8
$ws.Range("M98").Value = -867.8000000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 19762
$ws.Range("J86").Value = 19762
$ws.Range("L86").Value = 19762
$ws.Range("N86").Value = -22134
$ws.Range("H89").Value = 19762
$ws.Range("J89").Value = 19762
$ws.Range("L89").Value = 59286
$ws.Range("N89").Value = -71142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 490.6
$ws.Range("I16").Value = 490.6
$ws.Range("K16").Value = 490.6
$ws.Range("M16").Value = -320.6
$ws.Range("H22").Value = 1622
$ws.Range("I22").Value = 1496.1666
$ws.Range("K22").Value = 1496.1666
$ws.Range("M22").Value = -1201.1666
$ws.Range("H27").Value = 1622
$ws.Range("I27").Value = 1496.1666
$ws.Range("K27").Value = 1496.1666
$ws.Range("M27").Value = -1389.1666
$ws.Range("H122").Value = 10133.4
$ws.Range("I122").Value = 12354.182
$ws.Range("J122").Value = 4026.25
$ws.Range("K122").Value = 37062.546
$ws.Range("L122").Value = 12078.75
$ws.Range("M122").Value = -34612.546
$ws.Range("N122").Value = -16978.75
$ws.Range("H136").Value = 3144.6191
$ws.Range("I136").Value = 2478.9412
$ws.Range("J136").Value = 5973.75
$ws.Range("K136").Value = 7436.823600000001
$ws.Range("L136").Value = 17921.25
$ws.Range("M136").Value = -4886.823600000001
$ws.Range("N136").Value = -23021.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4307.357
$ws.Range("I81").Value = 4663.909
$ws.Range("K81").Value = 9327.817999999999
$ws.Range("M81").Value = -8266.817999999999
$ws.Range("H84").Value = 4307.357
$ws.Range("I84").Value = 4663.909
$ws.Range("K84").Value = 46639.09
$ws.Range("M84").Value = -41335.09
$ws.Range("H113").Value = 659
$ws.Range("I113").Value = 244.36363
$ws.Range("J113").Value = 1799.25
$ws.Range("K113").Value = 733.0908899999999
$ws.Range("L113").Value = 5397.75
$ws.Range("M113").Value = 1436.90911
$ws.Range("N113").Value = -9737.75
$ws.Range("H132").Value = 54964.406
$ws.Range("I132").Value = 50168.082
$ws.Range("K132").Value = 150504.246
$ws.Range("M132").Value = -147974.246
$ws.Range("H136").Value = 5815.778
$ws.Range("J136").Value = 1005
$ws.Range("L136").Value = 3015
$ws.Range("N136").Value = -8115
